$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.51199746131897
$ws.Range("E2").Value = 2079.833702892311
$ws.Range("F2").Value = 0.1208368762573507
$ws.Range("G2").Value = 0.1046978464327437
$ws.Range("H2").Value = 0.08670570476509028
$ws.Range("I2").Value = 0.07390905033000368
$ws.Range("J2").Value = 0.06582214971571611
$ws.Range("K2").Value = 0.05845205668232498
$ws.Range("L2").Value = 0.05262590202691655
$ws.Range("M2").Value = 0.05079474580444406
$ws.Range("N2").Value = 0.0476091867889111
$ws.Range("O2").Value = 0.04690039900655961
$ws.Range("P2").Value = 0.04569740408398586
$ws.Range("Q2").Value = 0.04446165239812615
$ws.Range("R2").Value = 0.04398749152540059
$ws.Range("S2").Value = 0.0425546719451282
$ws.Range("T2").Value = 0.04204789131131006
$ws.Range("U2").Value = 0.04129415094710601
$ws.Range("V2").Value = 0.04106983067256621
$ws.Range("W2").Value = 0.04078714127672181
$ws.Range("X2").Value = 0.04068608641518711
$ws.Range("Y2").Value = 0.0405425673078423

$ws.Range("C3").Value = 1.705002784729004
$ws.Range("E3").Value = 2038.937668118537
$ws.Range("F3").Value = 0.1208368762573507
$ws.Range("G3").Value = 0.1038994889086128
$ws.Range("H3").Value = 0.08680498896299618
$ws.Range("I3").Value = 0.07148908951482578
$ws.Range("J3").Value = 0.06053506220453538
$ws.Range("K3").Value = 0.05685773156125834
$ws.Range("L3").Value = 0.05345778748322966
$ws.Range("M3").Value = 0.05036450156899213
$ws.Range("N3").Value = 0.04858440223679898
$ws.Range("O3").Value = 0.04615648514149714
$ws.Range("P3").Value = 0.04463076432366069
$ws.Range("Q3").Value = 0.04356780593749897
$ws.Range("R3").Value = 0.04248362824236813
$ws.Range("S3").Value = 0.04162837034194226
$ws.Range("T3").Value = 0.04110084746464106
$ws.Range("U3").Value = 0.04080390334867321
$ws.Range("V3").Value = 0.04039534564531873
$ws.Range("W3").Value = 0.04006369357557517
$ws.Range("X3").Value = 0.03995632305817289
$ws.Range("Y3").Value = 0.03974537364753482

$ws.Range("C4").Value = 1.415000438690186
$ws.Range("E4").Value = 2052.174818160916
$ws.Range("F4").Value = 0.1208368762573507
$ws.Range("G4").Value = 0.1048996092386553
$ws.Range("H4").Value = 0.08557389090558423
$ws.Range("I4").Value = 0.07180160977952758
$ws.Range("J4").Value = 0.06533423509687646
$ws.Range("K4").Value = 0.05804220910539257
$ws.Range("L4").Value = 0.05381636819658082
$ws.Range("M4").Value = 0.05124276430580068
$ws.Range("N4").Value = 0.0488527778507403
$ws.Range("O4").Value = 0.04632522773310277
$ws.Range("P4").Value = 0.04423955746012088
$ws.Range("Q4").Value = 0.04394451853845266
$ws.Range("R4").Value = 0.04263941486137949
$ws.Range("S4").Value = 0.0414533729944305
$ws.Range("T4").Value = 0.04119520958578951
$ws.Range("U4").Value = 0.0409719498504208
$ws.Range("V4").Value = 0.04068907787013366
$ws.Range("W4").Value = 0.04042239214254753
$ws.Range("X4").Value = 0.04016467811261279
$ws.Range("Y4").Value = 0.04000340776142137

$ws.Range("C5").Value = 1.802998542785645
$ws.Range("E5").Value = 2008.763839129748
$ws.Range("F5").Value = 0.1208368762573507
$ws.Range("G5").Value = 0.1035703143463855
$ws.Range("H5").Value = 0.08292492226883169
$ws.Range("I5").Value = 0.07058879197472448
$ws.Range("J5").Value = 0.06435097485154793
$ws.Range("K5").Value = 0.05894780536718696
$ws.Range("L5").Value = 0.05480669207304058
$ws.Range("M5").Value = 0.0504595454847124
$ws.Range("N5").Value = 0.04861714865441696
$ws.Range("O5").Value = 0.04580839644987819
$ws.Range("P5").Value = 0.044181219632447
$ws.Range("Q5").Value = 0.04306761390738865
$ws.Range("R5").Value = 0.04166304960370325
$ws.Range("S5").Value = 0.04068396394608849
$ws.Range("T5").Value = 0.04017238656925051
$ws.Range("U5").Value = 0.03997955294821733
$ws.Range("V5").Value = 0.03972179487761587
$ws.Range("W5").Value = 0.03943950027157461
$ws.Range("X5").Value = 0.03929656056610993
$ws.Range("Y5").Value = 0.03915718984658377

$ws.Range("C6").Value = 1.445000648498535
$ws.Range("E6").Value = 2098.81305375008
$ws.Range("F6").Value = 0.1208368762573507
$ws.Range("G6").Value = 0.1040133364878942
$ws.Range("H6").Value = 0.08361693345235156
$ws.Range("I6").Value = 0.07182203448931962
$ws.Range("J6").Value = 0.06434495783869608
$ws.Range("K6").Value = 0.05973226035362898
$ws.Range("L6").Value = 0.05667219105066454
$ws.Range("M6").Value = 0.0533660064743588
$ws.Range("N6").Value = 0.04945652927863348
$ws.Range("O6").Value = 0.04705963957481313
$ws.Range("P6").Value = 0.04612284797348569
$ws.Range("Q6").Value = 0.04503372371564589
$ws.Range("R6").Value = 0.04366569357391826
$ws.Range("S6").Value = 0.04331279918777034
$ws.Range("T6").Value = 0.04227074418231552
$ws.Range("U6").Value = 0.04200742157649648
$ws.Range("V6").Value = 0.04165258842810322
$ws.Range("W6").Value = 0.04136708930821322
$ws.Range("X6").Value = 0.0409421428435264
$ws.Range("Y6").Value = 0.04091253516082027

$ws.Range("C7").Value = 1.54804253578186
$ws.Range("E7").Value = 2026.03092686823
$ws.Range("F7").Value = 0.1208368762573507
$ws.Range("G7").Value = 0.1047372430564227
$ws.Range("H7").Value = 0.08439128665721859
$ws.Range("I7").Value = 0.06988654839806971
$ws.Range("J7").Value = 0.06280280292160967
$ws.Range("K7").Value = 0.0580200278042592
$ws.Range("L7").Value = 0.05478901067073
$ws.Range("M7").Value = 0.05255540954034801
$ws.Range("N7").Value = 0.04944094479598955
$ws.Range("O7").Value = 0.04826787343748706
$ws.Range("P7").Value = 0.04583166538502016
$ws.Range("Q7").Value = 0.04446382736025691
$ws.Range("R7").Value = 0.04271577528226569
$ws.Range("S7").Value = 0.04187610314028494
$ws.Range("T7").Value = 0.04150959013153216
$ws.Range("U7").Value = 0.04085032024474697
$ws.Range("V7").Value = 0.04033571877496198
$ws.Range("W7").Value = 0.0398425065060766
$ws.Range("X7").Value = 0.03962578082754641
$ws.Range("Y7").Value = 0.03949378025084269

$ws.Range("C8").Value = 1.514000177383423
$ws.Range("E8").Value = 1974.683112871675
$ws.Range("F8").Value = 0.1208368762573507
$ws.Range("G8").Value = 0.1027741722434948
$ws.Range("H8").Value = 0.08158773619270858
$ws.Range("I8").Value = 0.06853498391267301
$ws.Range("J8").Value = 0.06030024660995889
$ws.Range("K8").Value = 0.0564400495448906
$ws.Range("L8").Value = 0.05178004705189682
$ws.Range("M8").Value = 0.04858837441721635
$ws.Range("N8").Value = 0.04575211372421718
$ws.Range("O8").Value = 0.04470062432103328
$ws.Range("P8").Value = 0.04322773246943647
$ws.Range("Q8").Value = 0.04189719580848435
$ws.Range("R8").Value = 0.04118804894329826
$ws.Range("S8").Value = 0.04018955622110286
$ws.Range("T8").Value = 0.03963089521685492
$ws.Range("U8").Value = 0.03926935852117994
$ws.Range("V8").Value = 0.03895193924283496
$ws.Range("W8").Value = 0.0388074620687009
$ws.Range("X8").Value = 0.03862381525560658
$ws.Range("Y8").Value = 0.03849284820412621

$ws.Range("C9").Value = 1.60501766204834
$ws.Range("E9").Value = 2075.667686173425
$ws.Range("F9").Value = 0.1208368762573507
$ws.Range("G9").Value = 0.103435246006345
$ws.Range("H9").Value = 0.08496726109185798
$ws.Range("I9").Value = 0.07313702618065152
$ws.Range("J9").Value = 0.06447139812128062
$ws.Range("K9").Value = 0.05814619266518838
$ws.Range("L9").Value = 0.05325352028123926
$ws.Range("M9").Value = 0.05062630443277006
$ws.Range("N9").Value = 0.04839848159449864
$ws.Range("O9").Value = 0.04647885980200781
$ws.Range("P9").Value = 0.04543066248641141
$ws.Range("Q9").Value = 0.04451101379934556
$ws.Range("R9").Value = 0.04310594117454055
$ws.Range("S9").Value = 0.04254996804049461
$ws.Range("T9").Value = 0.04210921269641679
$ws.Range("U9").Value = 0.04156394649014403
$ws.Range("V9").Value = 0.0412087068995173
$ws.Range("W9").Value = 0.04083577166733238
$ws.Range("X9").Value = 0.04046135840494005
$ws.Range("Y9").Value = 0.04046135840494005

$ws.Range("C10").Value = 1.515998125076294
$ws.Range("E10").Value = 2044.110134831099
$ws.Range("F10").Value = 0.1208368762573507
$ws.Range("G10").Value = 0.101467035841154
$ws.Range("H10").Value = 0.08456375833187919
$ws.Range("I10").Value = 0.07199023280576933
$ws.Range("J10").Value = 0.06429386321466224
$ws.Range("K10").Value = 0.05803216301735109
$ws.Range("L10").Value = 0.0526154203690779
$ws.Range("M10").Value = 0.05149959105390813
$ws.Range("N10").Value = 0.04888738908206169
$ws.Range("O10").Value = 0.04576861603422703
$ws.Range("P10").Value = 0.04380687031191129
$ws.Range("Q10").Value = 0.04304849710058537
$ws.Range("R10").Value = 0.04226810740608188
$ws.Range("S10").Value = 0.04158109062853486
$ws.Range("T10").Value = 0.04128014286313924
$ws.Range("U10").Value = 0.04078539784959424
$ws.Range("V10").Value = 0.0405537769722306
$ws.Range("W10").Value = 0.04011571137810088
$ws.Range("X10").Value = 0.03986602160099113
$ws.Range("Y10").Value = 0.03984620145869586

$ws.Range("C11").Value = 1.415998935699463
$ws.Range("E11").Value = 2011.177426735016
$ws.Range("F11").Value = 0.1208368762573507
$ws.Range("G11").Value = 0.1015972059999203
$ws.Range("H11").Value = 0.08102264183393888
$ws.Range("I11").Value = 0.06668870824595614
$ws.Range("J11").Value = 0.06308710562886734
$ws.Range("K11").Value = 0.05753844622481971
$ws.Range("L11").Value = 0.05417465232824378
$ws.Range("M11").Value = 0.04937579822456258
$ws.Range("N11").Value = 0.0477483613780106
$ws.Range("O11").Value = 0.04564941176363707
$ws.Range("P11").Value = 0.04393442406055611
$ws.Range("Q11").Value = 0.04199635390246673
$ws.Range("R11").Value = 0.04139168961131439
$ws.Range("S11").Value = 0.04115529381363594
$ws.Range("T11").Value = 0.04068040249330728
$ws.Range("U11").Value = 0.04035430737187006
$ws.Range("V11").Value = 0.03995719438863937
$ws.Range("W11").Value = 0.03967875427208662
$ws.Range("X11").Value = 0.03935213183677796
$ws.Range("Y11").Value = 0.03920423833791453
